$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I holds numbers-as-text in this sheet (inline strings); force text
# format so the COM layer doesn't silently convert the assigned values into
# numeric cells.
$ws.Range("I12:I21").NumberFormat = "@"

$ws.Range("A12").Value = 111378946
$ws.Range("I12").Value = "100"
$ws.Range("Q12").Value = 505602.791734456
$ws.Range("R12").Value = 6913005.013642685

$ws.Range("A13").Value = 111378856
$ws.Range("I13").Value = "10"
$ws.Range("Q13").Value = 505494.3524330241
$ws.Range("R13").Value = 6913043.848162009

$ws.Range("A14").Value = 111378884
$ws.Range("I14").Value = "50"
$ws.Range("Q14").Value = 505596.2310213979
$ws.Range("R14").Value = 6913034.263345711

$ws.Range("A15").Value = 111378874
$ws.Range("I15").Value = "50"
$ws.Range("Q15").Value = 505592.4968292552
$ws.Range("R15").Value = 6913042.152801346

$ws.Range("A16").Value = 111378866
$ws.Range("I16").Value = "10"
$ws.Range("Q16").Value = 505492.5216403615
$ws.Range("R16").Value = 6913025.731493607

$ws.Range("A17").Value = 111378913
$ws.Range("I17").Value = "25"
$ws.Range("Q17").Value = 505607.407264018
$ws.Range("R17").Value = 6913026.386327411

$ws.Range("A18").Value = 111378893
$ws.Range("I18").Value = "25"
$ws.Range("Q18").Value = 505612.5119866763
$ws.Range("R18").Value = 6913033.361683531

$ws.Range("A19").Value = 111378933
$ws.Range("I19").Value = "25"
$ws.Range("Q19").Value = 505597.6535686332
$ws.Range("R19").Value = 6913018.009825628

$ws.Range("A20").Value = 111378964
$ws.Range("I20").Value = "5"
$ws.Range("Q20").Value = 505627.1571942444
$ws.Range("R20").Value = 6912898.692122459

$ws.Range("A21").Value = 111378954
$ws.Range("I21").Value = "15"
$ws.Range("Q21").Value = 505590.6913760683
$ws.Range("R21").Value = 6913009.17353364
